$d = $word.ActiveDocument

# --- 1) "Total Hours: 16" -> "Total Hours: 18" -------------------------
$pTotal = $d.Paragraphs.Item(5).Range
$pTotal.Find.Execute("16", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "18", 2)

# --- 2) "Cumulative Hours: 16" -> "Cumulative Hours: 18" ----------------
$pCum = $d.Paragraphs.Item(6).Range
$pCum.Find.Execute("16", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "18", 2)

# --- 3) Merge the sentence split around the old bookmark location ------
$pStory = $d.Paragraphs.Item(8).Range
$pStory.Find.Execute("Her personal experience" + "s inspired me vastly and stirred my approach to work.", `
                      $true, $false, $false, $false, $false, `
                      $true, 1, $false, `
                      "Her personal experiences inspired me vastly and stirred my approach to work.", 2)

# --- 4) Move the "_GoBack" bookmark to sit right after the new "18" ----
#        in the Cumulative Hours paragraph (Word's last-edit marker).
$pCum2 = $d.Paragraphs.Item(6).Range
$pCum2.Find.Execute("18", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)
$goBack = $pCum2.Duplicate
$goBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBack)
